# Replace the bare word "него" in the phrase "в отношении него" with the
# gender-aware template placeholder "{{него_неё}}" (commit: "Add
# gender-based verb key obratlsya_as").
#
# Result: "...введении в отношении него процедуры..."
#      -> "...введении в отношении {{него_неё}} процедуры..."

$d = $word.ActiveDocument

# Anchor on the unique phrase "в отношении него" so we edit the right spot
# even if other occurrences of "него" exist elsewhere in the document.
$anchor = $d.Content
$found = $anchor.Find.Execute("в отношении него", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $phraseStart = $anchor.Start
    $phraseEnd = $anchor.End

    # "него" is the trailing 4 characters of the matched phrase.
    $wordRange = $d.Range($phraseEnd - 4, $phraseEnd)

    if ($wordRange.Text -eq "него") {
        $wordRange.Text = "{{него_неё}}"
    }
} else {
    # Fallback: the lead-in wording changed but the target word is still
    # present and unique - replace it directly.
    $anchor2 = $d.Content
    $found2 = $anchor2.Find.Execute("него", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $anchor2.Text = "{{него_неё}}"
    }
}
